# Sincronizando os dados entre documentos
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5 - B column text update
$ws.Range("B5").Value = "Cadastrar Conta ADM"

# Renumber CSU identifiers in column A (rows 15-31)
$ws.Range("A15").Value = "CSU-14"
$ws.Range("A16").Value = "CSU-15"
$ws.Range("A17").Value = "CSU-16"
$ws.Range("A18").Value = "CSU-17"
$ws.Range("A19").Value = "CSU-18"
$ws.Range("A20").Value = "CSU-19"
$ws.Range("A21").Value = "CSU-20"
$ws.Range("A22").Value = "CSU-21"
$ws.Range("A23").Value = "CSU-22"
$ws.Range("A24").Value = "CSU-23"
$ws.Range("A25").Value = "CSU-24"
$ws.Range("A26").Value = "CSU-25"
$ws.Range("A27").Value = "CSU-26"
$ws.Range("A28").Value = "CSU-27"
$ws.Range("A29").Value = "CSU-28"
$ws.Range("A30").Value = "CSU-29"
$ws.Range("A31").Value = "CSU-30"

# Fix casing of "Tipo de Ingresso" descriptions in column B (rows 24-27)
$ws.Range("B24").Value = "Cadastrar tipo de ingresso"
$ws.Range("B25").Value = "Exibir tipo de ingresso"
$ws.Range("B26").Value = "Editar tipo de ingresso"
$ws.Range("B27").Value = "Excluir tipo de ingresso"

# Update selection to match the saved view state
$ws.Range("B35").Select()
